$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new rows at the top of the data block (row 2),
#     pushing the existing data (rows 2-21) down to rows 4-23.
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Rows.Item(2).EntireRow.Insert()

# The inserted rows pick up formatting from the row above (the header),
# so strip that back out to match plain, unstyled data cells.
$ws.Range("A2:C3").ClearFormats()

# --- Step 2: populate the two newly inserted rows with their values.
$topData = @(
    @(0.2804546356201172, 0.4303635954856872, -0.691750168800354),
    @(0.1987819671630859, 0.2879692316055298, -0.9282988905906676)
)
for ($i = 0; $i -lt $topData.Length; $i++) {
    $row = 2 + $i
    $vals = $topData[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, 1 + $j).Value = $vals[$j]
    }
}

# --- Step 3: append eight brand-new data rows after the existing data
#     (which now ends at row 23).
$bottomData = @(
    @(-0.4514303207397461, -0.07753515243530271, -1.056098580360413),
    @(1.037992477416992, -1.273390769958496, 0.4362349510192871),
    @(0.0754270553588867, 1.646718859672546, 1.695090532302856),
    @(-0.2560558319091797, 0.3026316165924072, -0.4233262538909912),
    @(0.6335611343383789, 0.8106564879417419, -1.443797469139099),
    @(0.09285736083984369, 0.7357764840126038, -1.646607518196106),
    @(0.0882749557495117, 0.1726978719234466, -0.9354652166366576),
    @(0.2656211853027344, 0.4902379512786865, -0.8409426212310791)
)
for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $row = 24 + $i
    $vals = $bottomData[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, 1 + $j).Value = $vals[$j]
    }
}
